# Update Summary sheet (Sheet1) figures
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$expenses = $wb.Worksheets.Item("Expenses")

$summary.Range("A2").Value = 65000
$summary.Range("C2").Value = 4225
$summary.Range("D2").Value = 1990
$summary.Range("E2").Value = 2235

# Reorder / update Expenses sheet rows 3-6, and append new rows 7-9
$expenses.Range("A3").Value = "Car Insurance"
$expenses.Range("B3").Value = 180

$expenses.Range("A4").Value = "Groceries"
$expenses.Range("B4").Value = 300

$expenses.Range("A5").Value = "Utilities"
$expenses.Range("B5").Value = 0

$expenses.Range("A6").Value = "Savings"
$expenses.Range("B6").Value = 0

$expenses.Range("A7").Value = "Electricity"
$expenses.Range("B7").Value = 55

$expenses.Range("A8").Value = "Internet"
$expenses.Range("B8").Value = 55

$expenses.Range("A9").Value = "Miscellaneous"
$expenses.Range("B9").Value = 300
